$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.040887619778105
$ws.Range("D2").Value = 1.053408457646288
$ws.Range("E2").Value = 1.049754726995975
$ws.Range("F2").Value = 1.061874243139964
$ws.Range("I2").Value = 1.04432885826195
$ws.Range("J2").Value = 1.045971798341072
$ws.Range("K2").Value = 1.05615450815779
$ws.Range("L2").Value = 1.052510904218774
$ws.Range("M2").Value = 1.06459711976903
$ws.Range("N2").Value = 1.047457198269359

$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.041738859793997
$ws.Range("D3").Value = 1.05399211725576
$ws.Range("E3").Value = 1.050481294714198
$ws.Range("F3").Value = 1.062610504610459
$ws.Range("I3").Value = 1.04450753464971
$ws.Range("J3").Value = 1.046469330863596
$ws.Range("K3").Value = 1.05655179430724
$ws.Range("L3").Value = 1.053050002825333
$ws.Range("M3").Value = 1.065148284883627
$ws.Range("N3").Value = 1.047955437345133

$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.04229046237491
$ws.Range("D4").Value = 1.054370367610957
$ws.Range("E4").Value = 1.05095247351258
$ws.Range("F4").Value = 1.063087908635912
$ws.Range("I4").Value = 1.044622353916473
$ws.Range("J4").Value = 1.046791391578261
$ws.Range("K4").Value = 1.056808746313628
$ws.Range("L4").Value = 1.053399203888276
$ws.Range("M4").Value = 1.065505247151678
$ws.Range("N4").Value = 1.048277955422953

$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.042522544308871
$ws.Range("D5").Value = 1.054529521830665
$ws.Range("E5").Value = 1.051150804324366
$ws.Range("F5").Value = 1.063288845080097
$ws.Range("I5").Value = 1.044670432529011
$ws.Range("J5").Value = 1.046926814212198
$ws.Range("K5").Value = 1.056916739133787
$ws.Range("L5").Value = 1.053546094654594
$ws.Range("M5").Value = 1.0656553894196
$ws.Range("N5").Value = 1.048413570372564

$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.042561522873373
$ws.Range("D6").Value = 1.054556252547057
$ws.Range("E6").Value = 1.051184119382235
$ws.Range("F6").Value = 1.06332259695958
$ws.Range("I6").Value = 1.044678493893032
$ws.Range("J6").Value = 1.046949553889074
$ws.Range("K6").Value = 1.056934869815461
$ws.Range("L6").Value = 1.053570763297214
$ws.Range("M6").Value = 1.06568060333625
$ws.Range("N6").Value = 1.048436342342389

$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.042293562727693
$ws.Range("D7").Value = 1.054372493698526
$ws.Range("E7").Value = 1.050955122649541
$ws.Range("F7").Value = 1.063090592634635
$ws.Range("I7").Value = 1.044622997098308
$ws.Range("J7").Value = 1.046793200992584
$ws.Range("K7").Value = 1.056810189437827
$ws.Range("L7").Value = 1.05340116631135
$ws.Range("M7").Value = 1.065507253066097
$ws.Range("N7").Value = 1.048279767406851

$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.041175135015982
$ws.Range("D8").Value = 1.053605586151701
$ws.Range("E8").Value = 1.050000057261556
$ws.Range("F8").Value = 1.062122859163229
$ws.Range("I8").Value = 1.044389407200061
$ws.Range("J8").Value = 1.046139915317639
$ws.Range("K8").Value = 1.056288796615795
$ws.Range("L8").Value = 1.052693017966274
$ws.Range("M8").Value = 1.064783321017868
$ws.Range("N8").Value = 1.047625553991316

$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.039210475207469
$ws.Range("D9").Value = 1.05225875195054
$ws.Range("E9").Value = 1.048325168259044
$ws.Range("F9").Value = 1.060425291217736
$ws.Range("I9").Value = 1.043971727362314
$ws.Range("J9").Value = 1.044989751330764
$ws.Range("K9").Value = 1.055369188091922
$ws.Range("L9").Value = 1.051448057091855
$ws.Range("M9").Value = 1.063510196503324
$ws.Range("N9").Value = 1.046473756639654

$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.037904941814271
$ws.Range("D10").Value = 1.051364040638908
$ws.Range("E10").Value = 1.047214110017054
$ws.Range("F10").Value = 1.05929888229592
$ws.Range("I10").Value = 1.043689242821746
$ws.Range("J10").Value = 1.044223734501926
$ws.Range("K10").Value = 1.054755628712814
$ws.Range("L10").Value = 1.050620114620419
$ws.Range("M10").Value = 1.062663251171885
$ws.Range("N10").Value = 1.045706651979055

$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.037340658803419
$ws.Range("D11").Value = 1.050977399131944
$ws.Range("E11").Value = 1.046734346918927
$ws.Range("F11").Value = 1.058812419030221
$ws.Range("I11").Value = 1.043565977491487
$ws.Range("J11").Value = 1.043892237919165
$ws.Range("K11").Value = 1.054489852790369
$ws.Range("L11").Value = 1.050262107357847
$ws.Range("M11").Value = 1.062296963458942
$ws.Range("N11").Value = 1.045374684633123

$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.037131214197915
$ws.Range("D12").Value = 1.050833901607583
$ws.Range("E12").Value = 1.046556343741661
$ws.Range("F12").Value = 1.058631919258246
$ws.Range("I12").Value = 1.043520049663263
$ws.Range("J12").Value = 1.043769135818963
$ws.Range("K12").Value = 1.054391118079364
$ws.Range("L12").Value = 1.050129203784398
$ws.Range("M12").Value = 1.062160976484973
$ws.Range("N12").Value = 1.04525140771382

$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.037176133724158
$ws.Range("D13").Value = 1.050864676926136
$ws.Range("E13").Value = 1.04659451684521
$ws.Range("F13").Value = 1.058670628242981
$ws.Range("I13").Value = 1.043529907733325
$ws.Range("J13").Value = 1.043795540235408
$ws.Range("K13").Value = 1.054412297607136
$ws.Range("L13").Value = 1.050157708576034
$ws.Range("M13").Value = 1.06219014302124
$ws.Range("N13").Value = 1.045277849627565

$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.037323342870216
$ws.Range("D14").Value = 1.050965535156019
$ws.Range("E14").Value = 1.046719628984338
$ws.Range("F14").Value = 1.058797494885739
$ws.Range("I14").Value = 1.043562183966887
$ws.Range("J14").Value = 1.043882061631178
$ws.Range("K14").Value = 1.0544816916196
$ws.Range("L14").Value = 1.05025111994227
$ws.Range("M14").Value = 1.062285721324152
$ws.Range("N14").Value = 1.045364493893641

$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.037414063905472
$ws.Range("D15").Value = 1.051027692982074
$ws.Range("E15").Value = 1.046796741555953
$ws.Range("F15").Value = 1.058875687428429
$ws.Range("I15").Value = 1.043582051680096
$ws.Range("J15").Value = 1.043935374400266
$ws.Range("K15").Value = 1.05452444578818
$ws.Range("L15").Value = 1.050308683920112
$ws.Range("M15").Value = 1.062344619411071
$ws.Range("N15").Value = 1.045417882372976

$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.037942413091826
$ws.Range("D16").Value = 1.051389717248993
$ws.Range("E16").Value = 1.047245978558927
$ws.Range("F16").Value = 1.059331194396644
$ws.Range("I16").Value = 1.043697403636105
$ws.Range("J16").Value = 1.044245739038468
$ws.Range("K16").Value = 1.054773265386187
$ws.Range("L16").Value = 1.050643884993196
$ws.Range("M16").Value = 1.06268757000654
$ws.Range("N16").Value = 1.045728687764563

$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.038274107314048
$ws.Range("D17").Value = 1.051617014249545
$ws.Range("E17").Value = 1.047528131219208
$ws.Range("F17").Value = 1.059617265891757
$ws.Range("I17").Value = 1.04376950770449
$ws.Range("J17").Value = 1.044440475425419
$ws.Range("K17").Value = 1.054929317401514
$ws.Range("L17").Value = 1.05085428200508
$ws.Range("M17").Value = 1.062902814252112
$ws.Range("N17").Value = 1.04592370069952

$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.038467677390471
$ws.Range("D18").Value = 1.051749667328249
$ws.Range("E18").Value = 1.04769283451977
$ws.Range("F18").Value = 1.05978424983775
$ws.Range("I18").Value = 1.043811473299888
$ws.Range("J18").Value = 1.044554080467586
$ws.Range("K18").Value = 1.055020330093994
$ws.Range("L18").Value = 1.050977050899829
$ws.Range("M18").Value = 1.063028405420445
$ws.Range("N18").Value = 1.046037467073877

$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.038533696429582
$ws.Range("D19").Value = 1.051794911189756
$ws.Range("E19").Value = 1.047749015811205
$ws.Range("F19").Value = 1.059841207893398
$ws.Range("I19").Value = 1.043825766941836
$ws.Range("J19").Value = 1.044592819964439
$ws.Range("K19").Value = 1.055051361363236
$ws.Range("L19").Value = 1.051018919993917
$ws.Range("M19").Value = 1.063071235978069
$ws.Range("N19").Value = 1.046076261585259

$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.038238509467516
$ws.Range("D20").Value = 1.051592619703507
$ws.Range("E20").Value = 1.047497845610654
$ws.Range("F20").Value = 1.059586560360424
$ws.Range("I20").Value = 1.043761781077736
$ws.Range("J20").Value = 1.044419580119983
$ws.Range("K20").Value = 1.05491257550134
$ws.Range("L20").Value = 1.050831703433512
$ws.Range("M20").Value = 1.06287971613802
$ws.Range("N20").Value = 1.045902775720354

$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.03727998914662
$ws.Range("D21").Value = 1.05093583163814
$ws.Range("E21").Value = 1.046682780979872
$ws.Range("F21").Value = 1.05876013043936
$ws.Range("I21").Value = 1.04355268332305
$ws.Range("J21").Value = 1.043856582383838
$ws.Range("K21").Value = 1.054461257179283
$ws.Range("L21").Value = 1.050223610505409
$ws.Range("M21").Value = 1.062257573989026
$ws.Range("N21").Value = 1.045338978462846

$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.036678228356666
$ws.Range("D22").Value = 1.050523568618231
$ws.Range("E22").Value = 1.046171488351573
$ws.Range("F22").Value = 1.058241646630095
$ws.Range("I22").Value = 1.043420396420535
$ws.Range("J22").Value = 1.043502780454195
$ws.Range("K22").Value = 1.054177416878692
$ws.Range("L22").Value = 1.04984172035476
$ws.Range("M22").Value = 1.061866806219696
$ws.Range("N22").Value = 1.044984674093882

$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.036997147338476
$ws.Range("D23").Value = 1.050742051380978
$ws.Range("E23").Value = 1.046442422592719
$ws.Range("F23").Value = 1.058516397305848
$ws.Range("I23").Value = 1.043490601567025
$ws.Range("J23").Value = 1.043690320303881
$ws.Range("K23").Value = 1.054327892963912
$ws.Range("L23").Value = 1.050044125133972
$ws.Range("M23").Value = 1.0620739213423
$ws.Range("N23").Value = 1.045172480271666

$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.0382545942981
$ws.Range("D24").Value = 1.051603642318629
$ws.Range("E24").Value = 1.047511529976979
$ws.Range("F24").Value = 1.059600434487695
$ws.Range("I24").Value = 1.04376527269056
$ws.Range("J24").Value = 1.04442902175178
$ws.Range("K24").Value = 1.054920140475875
$ws.Range("L24").Value = 1.050841905570389
$ws.Range("M24").Value = 1.062890153049383
$ws.Range("N24").Value = 1.045912230760351

$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.039717647126326
$ws.Range("D25").Value = 1.052606389161361
$ws.Range("E25").Value = 1.048757200257622
$ws.Range("F25").Value = 1.060863227639441
$ws.Range("I25").Value = 1.044080421545827
$ws.Range("J25").Value = 1.0452869684067
$ws.Range("K25").Value = 1.055607020850229
$ws.Range("L25").Value = 1.051769557880132
$ws.Range("M25").Value = 1.063839018931781
$ws.Range("N25").Value = 1.046771395797929
